$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, pushing existing rows 106-122 down to 107-123.
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with the new weekly record.
$ws.Range("A106").Value = 1
$ws.Range("B106").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C106").Value = "Arica y Parinacota"
$ws.Range("D106").Value = 44985
$ws.Range("E106").Value = 15
$ws.Range("F106").Value = 100112038
$ws.Range("G106").Value = "Cebollín baby"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 300
$ws.Range("K106").Value = 1300
$ws.Range("L106").Value = 1500
$ws.Range("M106").Value = 1400
$ws.Range("N106").Value = '$/paquete 1,5 a 2 kilos'
$ws.Range("O106").Value = "Región de Arica y Parinacota"
$ws.Range("P106").Value = 700
$ws.Range("Q106").Value = 2
$ws.Range("R106").Value = "Hortaliza"
